# Commit: "adicionado utils, pathlib e status atrasado para não enviados"
#
# For every row whose "Situação" (column E) was left blank (i.e. the
# técnico never submitted the form), mark it as "Atrasado" (late) and
# highlight the cell with a new orange fill / white Arial text style
# (mirrors the other status styles already present in the sheet, e.g.
# the dark-green "Não enviado" style used in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column E that are currently empty ("sem situação") and must
# become "Atrasado".
$rows = @(4,6,9,11,12,13,14,15,16,18,19,22,23,24,25,26,27,31,32,38,39,40,41,42,43,45,47,48,49,51,53,54,55,58,60,61,62,63,64,67,68)

foreach ($r in $rows) {
    $cell = $ws.Range("E$r")

    # New "Atrasado" look: orange fill + bold-free white Arial text,
    # centered - consistent with the other status cells in the sheet.
    $cell.Font.Name = "Arial"
    $cell.Font.Color = 16777215        # white   (RGB 255,255,255 / hex FFFFFF)
    $cell.Interior.Color = 25855       # orange  (RGB 255,100,0   / hex FF6400)
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4108    # xlCenter

    $cell.Value = "Atrasado"
}
